$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity")

$ws.Range("B2").Value = "Opportunity Discussed Meeting"
$ws.Range("E2").Value = "Opportunity Discussed Meeting Description"

$ws.Range("E6").Select()
